$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 height --------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 58

# --- New styles are created in the SAME order the target workbook used them, so
# --- the generated fonts/fills/cellXfs line up with the authored styles.xml:
# ---   xf6 = red font (fontId 4)                -> used by C19
# ---   xf7 = existing "green" font + green fill -> used by B9, C9, C10
# ---   xf8 = default font + green fill          -> used by D9,E9,F9,D10,F10,F11
# ---   xf9 = existing "blue" font + green fill  -> used by E10, E11

# 1) Red font style (creates font 4 + cellXfs 6)
$ws.Range("C19").Value = 10024781
$ws.Range("C19").Font.Color = 255

# 2) Green-fill + existing "green" font (creates fill 2 + cellXfs 7)
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = 10003047
foreach ($addr in @("B9","C9","C10")) {
    $ws.Range($addr).Interior.Color = 5296274
}

# 3) Green-fill + default font (creates cellXfs 8)
foreach ($addr in @("D9","E9","F9","D10","F10","F11")) {
    $ws.Range($addr).Interior.Color = 5296274
}

# 4) Green-fill + existing "blue" font (creates cellXfs 9)
$ws.Range("C5").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
foreach ($addr in @("E10","E11")) {
    $ws.Range($addr).Interior.Color = 5296274
}

# --- Row 5 gains the value that used to sit in B12 (same plain "category" style)
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Value = 10003058

# --- Remove the now-empty source cells (value moved elsewhere) -----------------
$ws.Range("B8").Clear()
$ws.Range("B12").Clear()

# --- B13 gains the value that used to sit in B14 -------------------------------
$ws.Range("B13").Value = 10003851
$ws.Range("B14").Clear()

# --- B16 gains the value that used to sit in B15 -------------------------------
$ws.Range("B16").Value = 10012727
$ws.Range("B15").Clear()

# --- Old B18 (10024781) relocates into C19 (already written above); clear B18 -
$ws.Range("B18").Clear()

# --- C19 gains the D:F trio mirroring the row-17 pattern -----------------------
$ws.Range("D19").Value = 10015151
$ws.Range("E19").Value = 10014982
$ws.Range("F19").Value = 10040785

# --- New row 20 ------------------------------------------------------------------
$ws.Range("E20").Value = 10003041
$ws.Range("F20").Value = 10003057

# --- Selection / scroll position --------------------------------------------------
$ws.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
